$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data (rows 5-9), matching the existing table's column layout:
# A: SKU/code, B: DESCRIPCION, C: DEPARTAMENTO, D: Precio Vigente, E: Costo Promedio,
# F: P. OU TLET, G: COSTO OUTLET (formula), H: COSTO TOTAL

$ws.Range("A5").Value = "112345678"
$ws.Range("B5").Value = "MARIO"
$ws.Range("C5").Value = "LINEA BLANCA"
$ws.Range("D5").Value = 499990
$ws.Range("E5").Value = 289270
$ws.Range("F5").Value = 144635
$ws.Range("G5").Formula = "=F5*1.19"
$ws.Range("H5").Value = "50000"

$ws.Range("A6").Value = "114323454"
$ws.Range("B6").Value = "BOWSER"
$ws.Range("C6").Value = "LINEA BLANCA"
$ws.Range("D6").Value = 549990
$ws.Range("E6").Value = 358312
$ws.Range("F6").Value = 179156
$ws.Range("G6").Formula = "=F6*1.19"
$ws.Range("H6").Value = "399990"

$ws.Range("A7").Value = "123445234"
$ws.Range("B7").Value = "CONTROL"
$ws.Range("C7").Value = "LINEA BLANCA"
$ws.Range("D7").Value = 259990
$ws.Range("E7").Value = 163354
$ws.Range("F7").Value = 81677
$ws.Range("G7").Formula = "=F7*1.19"
$ws.Range("H7").Value = "1200"

$ws.Range("A8").Value = "124533535"
$ws.Range("B8").Value = "REPISA"
$ws.Range("C8").Value = "LINEA BLANCA"
$ws.Range("D8").Value = 499990
$ws.Range("E8").Value = 289270
$ws.Range("F8").Value = 144635
$ws.Range("G8").Formula = "=F8*1.19"
$ws.Range("H8").Value = "42000"

$ws.Range("A9").Value = "123232332"
$ws.Range("B9").Value = "REGLA"
$ws.Range("C9").Value = "LINEA BLANCA"
$ws.Range("D9").Value = 549990
$ws.Range("E9").Value = 358312
$ws.Range("F9").Value = 179156
$ws.Range("G9").Formula = "=F9*1.19"
$ws.Range("H9").Value = "30000"

# Update the active selection to mirror the author's final selection (row 10, full row)
$ws.Range("A10:XFD10").Select()
